{"js": "// Apply the \"Abstract Title\" / \"Abstract\" spacing / \"Footnote Block Text\"\n// style changes to the document's style sheet.\n\nconst styles = context.document.getStyles();\n\n// ---------------------------------------------------------------------\n// 1) New custom paragraph style \"Abstract Title\" (based on Normal, next\n//    paragraph style is \"Abstract\"): centered, kept-with-next heading-like\n//    run-in title that introduces an Abstract block.\n// ---------------------------------------------------------------------\ncontext.document.addStyle(\"Abstract Title\", \"Paragraph\");\nawait context.sync();\n\nconst abstractTitle = styles.getByNameOrNullObject(\"Abstract Title\");\nabstractTitle.load(\"nameLocal\");\nawait context.sync();\n\nabstractTitle.baseStyle = \"Normal\";\nabstractTitle.nextParagraphStyle = \"Abstract\";\nabstractTitle.quickStyle = true;\n\nabstractTitle.paragraphFormat.keepWithNext = true;\nabstractTitle.paragraphFormat.keepTogether = true;\nabstractTitle.paragraphFormat.alignment = \"Centered\";\nabstractTitle.paragraphFormat.spaceBefore = 15; // 300 twips\nabstractTitle.paragraphFormat.spaceAfter = 0;\n\nabstractTitle.font.size = 10;\nabstractTitle.font.bold = true;\nabstractTitle.font.color = \"#345A8A\";\n\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Existing \"Abstract\" style: tighten the space that precedes the\n//    paragraph (before=300 -> before=100 twips, i.e. 15pt -> 5pt); the\n//    trailing space is left untouched.\n// ---------------------------------------------------------------------\nconst abstractStyle = styles.getByNameOrNullObject(\"Abstract\");\nabstractStyle.load(\"nameLocal\");\nawait context.sync();\n\nabstractStyle.paragraphFormat.spaceBefore = 5; // 100 twips\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3) New \"Footnote Block Text\" style, based on / followed by the\n//    built-in \"Footnote Text\" style, with Block-Text-like indentation.\n// ---------------------------------------------------------------------\ncontext.document.addStyle(\"Footnote Block Text\", \"Paragraph\");\nawait context.sync();\n\nconst footnoteBlockText = styles.getByNameOrNullObject(\"Footnote Block Text\");\nfootnoteBlockText.load(\"nameLocal\");\nawait context.sync();\n\nfootnoteBlockText.baseStyle = \"Footnote Text\";\nfootnoteBlockText.nextParagraphStyle = \"Footnote Text\";\nfootnoteBlockText.priority = 9;\nfootnoteBlockText.unhideWhenUsed = true;\nfootnoteBlockText.quickStyle = true;\n\nfootnoteBlockText.paragraphFormat.spaceBefore = 5; // 100 twips\nfootnoteBlockText.paragraphFormat.spaceAfter = 5; // 100 twips\nfootnoteBlockText.paragraphFormat.firstLineIndent = 0;\nfootnoteBlockText.paragraphFormat.leftIndent = 24; // 480 twips\nfootnoteBlockText.paragraphFormat.rightIndent = 24; // 480 twips\n\nawait context.sync();\n", "ps1": "# Apply the \"Abstract Title\" / \"Abstract\" spacing / \"Footnote Block Text\"\n# style changes to the styles part of the active document.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) New custom paragraph style \"Abstract Title\" (based on Normal, next\n#    paragraph style is \"Abstract\"): centered, kept-with-next heading-like\n#    run-in title that introduces an Abstract block.\n# ---------------------------------------------------------------------\n$abstractTitle = $d.Styles.Add(\"Abstract Title\", 1)\n$abstractTitle.BaseStyle = \"Normal\"\n$abstractTitle.NextParagraphStyle = \"Abstract\"\n$abstractTitle.QuickStyle = $true\n\n$abstractTitle.ParagraphFormat.KeepWithNext = $true\n$abstractTitle.ParagraphFormat.KeepTogether = $true\n$abstractTitle.ParagraphFormat.Alignment = 1\n$abstractTitle.ParagraphFormat.SpaceBefore = 15\n$abstractTitle.ParagraphFormat.SpaceAfter = 0\n\n$abstractTitle.Font.Size = 10\n$abstractTitle.Font.SizeBi = 10\n$abstractTitle.Font.Bold = $true\n$abstractTitle.Font.Color = 9067060\n\n# ---------------------------------------------------------------------\n# 2) Existing \"Abstract\" style: tighten the space that precedes the\n#    paragraph (before=300 -> before=100 twips), leave the trailing\n#    space untouched.\n# ---------------------------------------------------------------------\n$abstract = $d.Styles(\"Abstract\")\n$abstract.ParagraphFormat.SpaceBefore = 5\n\n# ---------------------------------------------------------------------\n# 3) New \"Footnote Block Text\" style, based on / followed by the\n#    built-in \"Footnote Text\" style, with Block-Text-like indentation.\n# ---------------------------------------------------------------------\n$footnoteBlockText = $d.Styles.Add(\"Footnote Block Text\", 1)\n$footnoteBlockText.BaseStyle = \"Footnote Text\"\n$footnoteBlockText.NextParagraphStyle = \"Footnote Text\"\n$footnoteBlockText.Priority = 9\n$footnoteBlockText.UnhideWhenUsed = $true\n$footnoteBlockText.QuickStyle = $true\n\n$footnoteBlockText.ParagraphFormat.SpaceBefore = 5\n$footnoteBlockText.ParagraphFormat.SpaceAfter = 5\n$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0\n$footnoteBlockText.ParagraphFormat.LeftIndent = 24\n$footnoteBlockText.ParagraphFormat.RightIndent = 24\n"}
